$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F, rows 2:25 (bus voltage magnitude results; slack bus setpoint changed 1.05 -> 1.02 pu)
$blockBF = New-Object 'object[,]' 24,5
$blockBF[0,0] = 1.02
$blockBF[0,1] = 1.027032374748685
$blockBF[0,2] = 1.030510707607442
$blockBF[0,3] = 1.030662890396634
$blockBF[0,4] = 1.036818059722732
$blockBF[1,0] = 1.02
$blockBF[1,1] = 1.028411716140358
$blockBF[1,2] = 1.031522547782159
$blockBF[1,3] = 1.031990603606696
$blockBF[1,4] = 1.03834280176057
$blockBF[2,0] = 1.02
$blockBF[2,1] = 1.029302875941649
$blockBF[2,2] = 1.032175871949847
$blockBF[2,3] = 1.032848746489432
$blockBF[2,4] = 1.0393281412644
$blockBF[3,0] = 1.02
$blockBF[3,1] = 1.029677198577296
$blockBF[3,2] = 1.032450196806962
$blockBF[3,3] = 1.033209281246208
$blockBF[3,4] = 1.039742080887668
$blockBF[4,0] = 1.02
$blockBF[4,1] = 1.029740030352725
$blockBF[4,2] = 1.032496237754067
$blockBF[4,3] = 1.033269803369626
$blockBF[4,4] = 1.039811565954142
$blockBF[5,0] = 1.02
$blockBF[5,1] = 1.029307878912507
$blockBF[5,2] = 1.032179538794666
$blockBF[5,3] = 1.032853564863924
$blockBF[5,4] = 1.039333673504584
$blockBF[6,0] = 1.02
$blockBF[6,1] = 1.027498816248034
$blockBF[6,2] = 1.030852956583146
$blockBF[6,3] = 1.031111802607375
$blockBF[6,4] = 1.037333620775637
$blockBF[7,0] = 1.02
$blockBF[7,1] = 1.024300275656996
$blockBF[7,2] = 1.028504436737592
$blockBF[7,3] = 1.028034882009094
$blockBF[7,4] = 1.033799223562198
$blockBF[8,0] = 1.02
$blockBF[8,1] = 1.022160304355713
$blockBF[8,2] = 1.026931194708914
$blockBF[8,3] = 1.025978084206201
$blockBF[8,4] = 1.031435755634544
$blockBF[9,0] = 1.02
$blockBF[9,1] = 1.021231780648191
$blockBF[9,2] = 1.026248118883846
$blockBF[9,3] = 1.025086084822186
$blockBF[9,4] = 1.030410539283866
$blockBF[10,0] = 1.02
$blockBF[10,1] = 1.020886592234334
$blockBF[10,2] = 1.025994111417604
$blockBF[10,3] = 1.024754540830575
$blockBF[10,4] = 1.030029446259404
$blockBF[11,0] = 1.02
$blockBF[11,1] = 1.020960649719135
$blockBF[11,2] = 1.026048609746468
$blockBF[11,3] = 1.024825668008498
$blockBF[11,4] = 1.030111204909697
$blockBF[12,0] = 1.02
$blockBF[12,1] = 1.021203253295035
$blockBF[12,2] = 1.026227128340455
$blockBF[12,3] = 1.025058683730599
$blockBF[12,4] = 1.030379043819464
$blockBF[13,0] = 1.02
$blockBF[13,1] = 1.021352690262525
$blockBF[13,2] = 1.026337081943357
$blockBF[13,3] = 1.025202223628722
$blockBF[13,4] = 1.030544030578657
$blockBF[14,0] = 1.02
$blockBF[14,1] = 1.022221886694671
$blockBF[14,2] = 1.026976488806449
$blockBF[14,3] = 1.026037253377013
$blockBF[14,4] = 1.031503756797426
$blockBF[15,0] = 1.02
$blockBF[15,1] = 1.022766596074864
$blockBF[15,2] = 1.027377072935709
$blockBF[15,3] = 1.026560668116827
$blockBF[15,4] = 1.032105274528832
$blockBF[16,0] = 1.02
$blockBF[16,1] = 1.023084132841683
$blockBF[16,2] = 1.027610548531292
$blockBF[16,3] = 1.026865833180001
$blockBF[16,4] = 1.032455954717001
$blockBF[17,0] = 1.02
$blockBF[17,1] = 1.02319237386256
$blockBF[17,2] = 1.027690127559844
$blockBF[17,3] = 1.026969864064809
$blockBF[17,4] = 1.032575498190441
$blockBF[18,0] = 1.02
$blockBF[18,1] = 1.022708172897072
$blockBF[18,2] = 1.027334112529959
$blockBF[18,3] = 1.026504524565891
$blockBF[18,4] = 1.032040755525749
$blockBF[19,0] = 1.02
$blockBF[19,1] = 1.021131820748424
$blockBF[19,2] = 1.02617456690401
$blockBF[19,3] = 1.024990072415297
$blockBF[19,4] = 1.030300179791851
$blockBF[20,0] = 1.02
$blockBF[20,1] = 1.020139004448376
$blockBF[20,2] = 1.02544387672295
$blockBF[20,3] = 1.024036624976011
$blockBF[20,4] = 1.02920417467241
$blockBF[21,0] = 1.02
$blockBF[21,1] = 1.020665478957207
$blockBF[21,2] = 1.025831386173608
$blockBF[21,3] = 1.024542186297556
$blockBF[21,4] = 1.029785345731841
$blockBF[22,0] = 1.02
$blockBF[22,1] = 1.022734572382046
$blockBF[22,2] = 1.027353525039497
$blockBF[22,3] = 1.026529893833936
$blockBF[22,4] = 1.032069909426837
$blockBF[23,0] = 1.02
$blockBF[23,1] = 1.02512848783676
$blockBF[23,2] = 1.029112901643891
$blockBF[23,3] = 1.028831287747729
$blockBF[23,4] = 1.034714186596857
$ws.Range("B2:F25").Value = $blockBF

# Columns I:N, rows 2:25
$blockIN = New-Object 'object[,]' 24,6
$blockIN[0,0] = 1.032200491400207
$blockIN[0,1] = 1.032192506076122
$blockIN[0,2] = 1.033321508912322
$blockIN[0,3] = 1.033473251036849
$blockIN[0,4] = 1.039610715061926
$blockIN[0,5] = 1.014615771517216
$blockIN[1,0] = 1.032550126842537
$blockIN[1,1] = 1.033209651014153
$blockIN[1,2] = 1.034141172593327
$blockIN[1,3] = 1.034607972194407
$blockIN[1,4] = 1.040943242613162
$blockIN[1,5] = 1.014964502449705
$blockIN[2,0] = 1.032774040350305
$blockIN[2,1] = 1.033866050625694
$blockIN[2,2] = 1.034669546099249
$blockIN[2,3] = 1.035340706593121
$blockIN[2,4] = 1.041803718468103
$blockIN[2,5] = 1.015189253489353
$blockIN[3,0] = 1.032867618955144
$blockIN[3,1] = 1.034141583755264
$blockIN[3,2] = 1.034891197992805
$blockIN[3,3] = 1.03564839228731
$blockIN[3,4] = 1.042165047789203
$blockIN[3,5] = 1.015283524492885
$blockIN[4,0] = 1.032883298740819
$blockIN[4,1] = 1.0341878226376
$blockIN[4,2] = 1.034928386486132
$blockIN[4,3] = 1.035700033338461
$blockIN[4,4] = 1.04222569243398
$blockIN[4,5] = 1.015299340481477
$blockIN[5,0] = 1.03277529292941
$blockIN[5,1] = 1.033869733947049
$blockIN[5,2] = 1.03467250969021
$blockIN[5,3] = 1.035344819296848
$blockIN[5,4] = 1.041808548190026
$blockIN[5,5] = 1.015190513983552
$blockIN[6,0] = 1.032319134585849
$blockIN[6,1] = 1.03253662226804
$blockIN[6,2] = 1.03359893500338
$blockIN[6,3] = 1.033857050048845
$blockIN[6,4] = 1.040061416901018
$blockIN[6,5] = 1.014733814220308
$blockIN[7,0] = 1.031497447562037
$blockIN[7,1] = 1.030173824139596
$blockIN[7,2] = 1.031691661461138
$blockIN[7,3] = 1.031223663848483
$blockIN[7,4] = 1.036968997308809
$blockIN[7,5] = 1.013922085096372
$blockIN[8,0] = 1.030937524353332
$blockIN[8,1] = 1.028589147148303
$blockIN[8,2] = 1.030409511081577
$blockIN[8,3] = 1.029459876251819
$blockIN[8,4] = 1.034897741372589
$blockIN[8,5] = 1.013376158829105
$blockIN[9,0] = 1.03069216842604
$blockIN[9,1] = 1.027900653986625
$blockIN[9,2] = 1.029851754805388
$blockIN[9,3] = 1.028694125296
$blockIN[9,4] = 1.033998488339301
$blockIN[9,5] = 1.013138614152689
$blockIN[10,0] = 1.030600593516523
$blockIN[10,1] = 1.027644563310645
$blockIN[10,2] = 1.029644187882219
$blockIN[10,3] = 1.028409382029356
$blockIN[10,4] = 1.033664099542775
$blockIN[10,5] = 1.013050204161133
$blockIN[11,0] = 1.03062025653863
$blockIN[11,1] = 1.027699511722198
$blockIN[11,2] = 1.029688729473152
$blockIN[11,3] = 1.028470474516709
$blockIN[11,4] = 1.033735843770758
$blockIN[11,5] = 1.013069176367132
$blockIN[12,0] = 1.03068460777997
$blockIN[12,1] = 1.027879492697848
$blockIN[12,2] = 1.029834605268322
$blockIN[12,3] = 1.028670594671951
$blockIN[12,4] = 1.033970855172268
$blockIN[12,5] = 1.013131309742169
$blockIN[13,0] = 1.030724198490659
$blockIN[13,1] = 1.027990337889125
$blockIN[13,2] = 1.029924432165906
$blockIN[13,3] = 1.028793854146496
$blockIN[13,4] = 1.034115604696799
$blockIN[13,5] = 1.013169568879321
$blockIN[14,0] = 1.030953746409107
$blockIN[14,1] = 1.028634790877418
$blockIN[14,2] = 1.03044647287549
$blockIN[14,3] = 1.029510653574986
$blockIN[14,4] = 1.034957370911494
$blockIN[14,5] = 1.013391899399435
$blockIN[15,0] = 1.031096956272736
$blockIN[15,1] = 1.029038414882822
$blockIN[15,2] = 1.030773242349507
$blockIN[15,3] = 1.029959738153864
$blockIN[15,4] = 1.035484744467644
$blockIN[15,5] = 1.013531050970865
$blockIN[16,0] = 1.031180207981838
$blockIN[16,1] = 1.029273618697437
$blockIN[16,2] = 1.030963593133319
$blockIN[16,3] = 1.030221487160043
$blockIN[16,4] = 1.035792122899375
$blockIN[16,5] = 1.013612104409864
$blockIN[17,0] = 1.031208547207957
$blockIN[17,1] = 1.029353779419388
$blockIN[17,2] = 1.031028455840631
$blockIN[17,3] = 1.030310704005317
$blockIN[17,4] = 1.035896892284229
$blockIN[17,5] = 1.013639722710427
$blockIN[18,0] = 1.031081620206702
$blockIN[18,1] = 1.028995132977254
$blockIN[18,2] = 1.030738208793672
$blockIN[18,3] = 1.029911575782274
$blockIN[18,4] = 1.035428186085338
$blockIN[18,5] = 1.013516132851391
$blockIN[19,0] = 1.03066567007457
$blockIN[19,1] = 1.027826502593365
$blockIN[19,2] = 1.029791659321917
$blockIN[19,3] = 1.028611672858923
$blockIN[19,4] = 1.033901660343196
$blockIN[19,5] = 1.013113017866103
$blockIN[20,0] = 1.030401605847376
$blockIN[20,1] = 1.027089689373579
$blockIN[20,2] = 1.029194259381184
$blockIN[20,3] = 1.027792580023609
$blockIN[20,4] = 1.032939750514223
$blockIN[20,5] = 1.012858548215208
$blockIN[21,0] = 1.030541832807074
$blockIN[21,1] = 1.027480484054002
$blockIN[21,2] = 1.029511168693153
$blockIN[21,3] = 1.028226968499116
$blockIN[21,4] = 1.033449880916403
$blockIN[21,5] = 1.012993544230197
$blockIN[22,0] = 1.031088550780682
$blockIN[22,1] = 1.02901469089737
$blockIN[22,2] = 1.030754039717083
$blockIN[22,3] = 1.029933338885917
$blockIN[22,4] = 1.035453743093601
$blockIN[22,5] = 1.013522874051507
$blockIN[23,0] = 1.031712003397788
$blockIN[23,1] = 1.03078631423571
$blockIN[23,2] = 1.032186595948527
$blockIN[23,3] = 1.031905879980155
$blockIN[23,4] = 1.037770132197861
$blockIN[23,5] = 1.014132771245051
$ws.Range("I2:N25").Value = $blockIN
